# "better logic and add granary item"
# - bump effect.propertyType (F2) and effect.value (G2) for the lucky_potion row
# - move the active selection to H10
# - resize the workbook window (best effort; host may not persist window chrome)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Data edits: row 2 is the "lucky_potion" expendable item.
#   F2 = effect.propertyType : 8  -> 10
#   G2 = effect.value        : 10 -> 150
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 150

# Move the selection/active cell to H10
$ws.Range("H10").Select()

# Resize the workbook window (mirrors windowHeight going from 13160 to 16760
# in the underlying bookViews entry)
$excel.ActiveWindow.Height = 16760
